# Updates the cryptos list: refreshes per-coin price (column D) and
# 1h volume/change percentage (column E) figures, and also reflects a
# shift of several ranking rows (48-51) to new coins, matching the
# upstream GitHub Actions data refresh.
#
# Note: some price strings look like plain decimal numbers (e.g. "19.92");
# Excel would otherwise auto-convert them to numeric cells. We force the
# cell to Text format immediately before assigning the value, then reset
# the cell style back to Normal so no stray formatting remains, keeping
# the cell's effective style identical to before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.206.27'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '1.643.59'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +1.14%  '
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0848'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '1.873.77'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '1.625.58'
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.15'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("E15").Value = '  +3.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("D17").Value = '27.194.53'
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '218.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("E21").Value = '  +3.03%  '
$ws.Range("E22").Value = '  +4.91%  '
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.73'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0508'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("E33").Value = '  +1.18%  '
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("D35").Value = '1.261.43'
$ws.Range("E35").Value = '  +1.56%  '
$ws.Range("E36").Value = '  +0.75%  '
$ws.Range("E37").Value = '  +2.28%  '
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("E39").Value = '  +2.03%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.809'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("E42").Value = '  +6.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").Value = '1.783.69'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0514'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.46%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0975'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.406'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.03%  '
